$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "Save" column in H1, matching the style of the
# neighboring header cells (e.g. G1: bold, bordered, centered).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Fill H2:H6 with the new "Save" column values (all 0 for now)
$ws.Range("H2:H6").Value = 0
